$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$newHeaders = @(
    "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404",
    "Segment ID_FV2404", "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404", "Bedingung_FV2404", "diff",
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
    "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# 2. Turn the used range into an Excel Table ("Table1") so headers are exposed as table columns
$fullRange = $ws.Range("A1:U85")
$tbl = $ws.ListObjects.Add(1, $fullRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split below row 1, keep top-left cell at A2)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
